# Leave Card update for MAMARIL JOSEFINA PEREY
# - Add a new SL (Sick Leave) entry for 3/25,26/2023 on the existing
#   March 2023 continuation row (row 98).
# - Extend the monthly PERIOD date column (A) through May 2025 (rows 99-124).
# - Grow Table1 by one row (A8:K133 -> A8:K134), preserving the special
#   "final row" border formatting on the new bottom row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- 1) Record the new sick-leave entry on row 98 -------------------------
$ws.Range("B98").Value = "SL(2-0-0)"
$ws.Range("H98").Value = 2
$ws.Range("K98").Value = "3/25,26/2023"

# --- 2) Fill in the PERIOD (month) dates for rows 99-124 -------------------
$monthDates = @(45017,45047,45078,45108,45139,45170,45200,45231,45261,45292,45323,45352,45383,45413,45444,45474,45505,45536,45566,45597,45627,45658,45689,45717,45748,45778)
$row = 99
foreach ($d in $monthDates) {
    $ws.Cells.Item($row, 1).Value = $d
    $row = $row + 1
}

# --- 3) Grow the table by one row, keeping the bottom row's distinct style -
# Move the current last row (133, with its special bottom-border styling)
# down to row 134 ...
$ws.Range("A133:K133").Copy($ws.Range("A134:K134"))
# ... then restyle row 133 like a normal interior row (copy row 132's look).
$ws.Range("A132:K132").Copy($ws.Range("A133:K133"))

# Resize the table to officially include the new row.
$lo.Resize($ws.Range("A8:K134"))

# Restore the calculated-column formula on the new bottom row (G134), which
# the row-132 style copy did not carry a formula for.
$ws.Range("G134").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

$excel.CalculateFullRebuild()
